# Generate Report for Handoff
# Updates the "Status" / handback-sync text to "Ready for handoff", refreshes the
# handback timestamps, and records "version not latest" error details for the
# two e2e files on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$overviewTimestamp = "2016-09-05 11:45:47"
$zhHandbackTimestamp = "2016-09-05 11:45:36"
$deHandbackTimestamp = "2016-09-05 11:45:47"

$err5b93 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb780f092c522aa03698e494d6dfcb8754686810/e2e/5b935408-90fc-4b55-a235-4bf052352988.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2ba69d0acbb3e9ff0893b75a35ab257e80d0bee/e2e/5b935408-90fc-4b55-a235-4bf052352988.md."
$errfdbb = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb780f092c522aa03698e494d6dfcb8754686810/e2e/fdbbaee0-fc80-46b6-902e-aaf59d871475.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2ba69d0acbb3e9ff0893b75a35ab257e80d0bee/e2e/fdbbaee0-fc80-46b6-902e-aaf59d871475.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $overviewTimestamp
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewTimestamp

$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("H2").Value = $zhHandbackTimestamp
$wsZh.Range("P2").Value = $err5b93
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("H3").Value = $zhHandbackTimestamp
$wsZh.Range("P3").Value = $errfdbb

$wsZh.Columns.Item(3).ColumnWidth = 16.25
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("H2").Value = $deHandbackTimestamp
$wsDe.Range("P2").Value = $err5b93
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("H3").Value = $deHandbackTimestamp
$wsDe.Range("P3").Value = $errfdbb

$wsDe.Columns.Item(3).ColumnWidth = 16.25
$wsDe.Columns.Item(16).ColumnWidth = 39.17
